$d = $word.ActiveDocument

# The document ends with a paragraph ("Adding line 3") that carries the
# Word-managed "_GoBack" bookmark right after its text. We want to add a new
# blank paragraph followed by a paragraph containing "Adding line 4", with
# the "_GoBack" bookmark now wrapping the newly typed text (mirroring how
# Word itself relocates _GoBack to the most recently edited text).

# Detach the existing bookmark so it doesn't interfere with the paragraph
# split below (if it isn't present for some reason, just continue).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$newText = "Adding line 4"

# Move to the very end of the document and type: <enter><enter>Adding line 4
$sel = $word.Selection
$sel.EndKey(6) | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeParagraph() | Out-Null
$sel.TypeText($newText) | Out-Null

# Re-create the "_GoBack" bookmark around the freshly typed text.
$lastPara = $d.Paragraphs.Last
$bmStart = $lastPara.Range.Start
$bmEnd = $bmStart + $newText.Length
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
